$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.364.04"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "2.640.94"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D5").Value = "'599.82"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").Value = "'154.77"
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "2.639.92"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "'0.137"
$ws.Range("E10").Value = "  +7.06%  "
$ws.Range("E11").Value = "  -0.53%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'28.08"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "'0.0000187"
$ws.Range("E15").Value = "  +3.38%  "
$ws.Range("D16").Value = "3.120.21"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "68.287.72"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "2.636.80"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "'11.46"
$ws.Range("E19").Value = "  +3.96%  "
$ws.Range("D20").Value = "'366.97"
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'7.44"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("E22").Value = "  +2.50%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("D25").Value = "'73.57"
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("E27").Value = "  +1.24%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.770.21"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0000105"
$ws.Range("E29").Value = "  +6.35%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'573.85"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "'1.43"
$ws.Range("E32").Value = "  +5.10%  "
$ws.Range("E34").Value = "  +2.72%  "
$ws.Range("E35").Value = "  +3.35%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").Value = "'160.57"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "'19.34"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "'5.41"
$ws.Range("E42").Value = "  +3.75%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  +4.01%  "
$ws.Range("E44").Value = "  +3.65%  "
$ws.Range("E45").Value = "  +14.09%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'40.52"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("D48").Value = "'158.43"
$ws.Range("E48").Value = "  +3.57%  "
$ws.Range("E49").Value = "  +3.30%  "
$ws.Range("D50").Value = "'1.71"
$ws.Range("E50").Value = "  +2.55%  "
$ws.Range("D51").Value = "'22.03"
$ws.Range("E51").Value = "  +3.40%  "
